{"js": "// Update the date paragraph and every populated table cell to the new\n// values used in the target revision.\n//\n// The replacement is driven by the exact *old* text each paragraph/cell\n// currently holds (all old values are unique in this document), and every\n// search is scoped to its own paragraph/cell so that values that are both\n// an old value somewhere AND a new value somewhere else (e.g. \"65\u00f72=32, 1\"\n// and \"10\u00f77=1, 3\" each appear once as an old value and once as a new\n// value) can never cross-match each other. insertText(..., \"Replace\") on\n// the matched range only swaps the run's text, leaving the run/paragraph\n// formatting (fonts, size, alignment) untouched.\n\nconst dateUpdate = { oldText: \"2024-11-16 Saturday\", newText: \"2024-11-17 Sunday\" };\n\n// Row-major list of [oldText, newText] pairs for the populated table cells,\n// in document (reading) order.\nconst cellUpdates = [\n  [\"65\u00f78=8, 1\", \"74\u00f78=9, 2\"],\n  [\"95\u00f75=19, 0\", \"94\u00f78=11, 6\"],\n  [\"80\u00f79=8, 8\", \"65\u00f72=32, 1\"],\n  [\"80\u00f76=13, 2\", \"94\u00f76=15, 4\"],\n  [\"10\u00f77=1, 3\", \"25\u00f73=8, 1\"],\n\n  [\"41\u00f75=8, 1\", \"54\u00f77=7, 5\"],\n  [\"71\u00f77=10, 1\", \"89\u00f76=14, 5\"],\n  [\"70\u00f79=7, 7\", \"74\u00f75=14, 4\"],\n  [\"10\u00f72=5, 0\", \"87\u00f73=29, 0\"],\n  [\"60\u00f72=30, 0\", \"91\u00f76=15, 1\"],\n\n  [\"94\u00f74=23, 2\", \"61\u00f76=10, 1\"],\n  [\"35\u00f79=3, 8\", \"94\u00f75=18, 4\"],\n  [\"46\u00f74=11, 2\", \"41\u00f76=6, 5\"],\n  [\"38\u00f77=5, 3\", \"10\u00f77=1, 3\"],\n  [\"65\u00f72=32, 1\", \"78\u00f73=26, 0\"],\n\n  [\"90\u00f73=30, 0\", \"27\u00f75=5, 2\"],\n  [\"45\u00f77=6, 3\", \"13\u00f75=2, 3\"],\n  [\"39\u00f77=5, 4\", \"79\u00f78=9, 7\"],\n  [\"14\u00f79=1, 5\", \"54\u00f73=18, 0\"],\n  [\"93\u00f79=10, 3\", \"15\u00f78=1, 7\"],\n\n  [\"17\u00f72=8, 1\", \"70\u00f75=14, 0\"],\n  [\"39\u00f74=9, 3\", \"41\u00f78=5, 1\"],\n  [\"36\u00f78=4, 4\", \"73\u00f76=12, 1\"],\n  [\"68\u00f73=22, 2\", \"41\u00f77=5, 6\"],\n  [\"88\u00f76=14, 4\", \"89\u00f74=22, 1\"],\n];\n\nasync function replaceInScope(scope, oldText, newText) {\n  const results = scope.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${JSON.stringify(oldText)}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\n// 1) Update the date heading (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nawait replaceInScope(paragraphs.items[0], dateUpdate.oldText, dateUpdate.newText);\nawait context.sync();\n\n// 2) Update each populated cell of the (single) table, walking rows/\n//    columns in order and skipping the blank spacer rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nlet updateIndex = 0;\nfor (let r = 0; r < table.rowCount && updateIndex < cellUpdates.length; r++) {\n  const cell0 = table.getCellOrNullObject(r, 0);\n  cell0.load([\"isNullObject\"]);\n  cell0.body.load(\"text\");\n  await context.sync();\n  if (cell0.isNullObject || !cell0.body.text || cell0.body.text.trim().length === 0) {\n    // Blank spacer row - nothing to update here.\n    continue;\n  }\n\n  for (let c = 0; c < 5; c++) {\n    const [oldText, newText] = cellUpdates[updateIndex];\n    const cell = table.getCellOrNullObject(r, c);\n    cell.load(\"isNullObject\");\n    await context.sync();\n    if (cell.isNullObject) {\n      throw new Error(`Expected a cell at row ${r}, col ${c}`);\n    }\n    await replaceInScope(cell.body, oldText, newText);\n    updateIndex++;\n  }\n}\n\nawait context.sync();\n\nif (updateIndex !== cellUpdates.length) {\n  throw new Error(`Only updated ${updateIndex} of ${cellUpdates.length} cells`);\n}\n", "ps1": "# Update the date paragraph and every populated table cell to the new\n# values used in the target revision.\n#\n# Unlike a find/replace, each cell is addressed by its absolute (row, col)\n# position in the table, so there is no risk of a new value (e.g.\n# \"65\u00f72=32, 1\") being mistaken for old text still waiting to be replaced\n# elsewhere in the document (it is also the OLD value of a different\n# cell). Setting `Cell.Range.Text` / `Paragraph.Range.Text` only rewrites\n# the run's text and leaves existing run/paragraph formatting (fonts,\n# size, alignment) untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading (first paragraph of the body). Note: no\n#    trailing `r here - that would insert an extra paragraph break and\n#    leave a stray empty paragraph behind.\n$d.Paragraphs.Item(1).Range.Text = \"2024-11-17 Sunday\"\n\n# 2) Update each populated cell of the (single) table, walking rows/\n#    columns in order and skipping the blank spacer rows.\n$newValues = @(\n  @(\"74\u00f78=9, 2\", \"94\u00f78=11, 6\", \"65\u00f72=32, 1\", \"94\u00f76=15, 4\", \"25\u00f73=8, 1\"),\n  @(\"54\u00f77=7, 5\", \"89\u00f76=14, 5\", \"74\u00f75=14, 4\", \"87\u00f73=29, 0\", \"91\u00f76=15, 1\"),\n  @(\"61\u00f76=10, 1\", \"94\u00f75=18, 4\", \"41\u00f76=6, 5\", \"10\u00f77=1, 3\", \"78\u00f73=26, 0\"),\n  @(\"27\u00f75=5, 2\", \"13\u00f75=2, 3\", \"79\u00f78=9, 7\", \"54\u00f73=18, 0\", \"15\u00f78=1, 7\"),\n  @(\"70\u00f75=14, 0\", \"41\u00f78=5, 1\", \"73\u00f76=12, 1\", \"41\u00f77=5, 6\", \"89\u00f74=22, 1\")\n)\n\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$dataRowIdx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  $firstCellText = $tbl.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n  if ($firstCellText.Length -eq 0) {\n    continue\n  }\n  if ($dataRowIdx -ge $newValues.Count) {\n    continue\n  }\n  $rowValues = $newValues[$dataRowIdx]\n  for ($c = 1; $c -le $colCount; $c++) {\n    $tbl.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n  }\n  $dataRowIdx++\n}\n"}
